$d = $word.ActiveDocument

# Questão 2: merge two runs into "Questão 2 =  falsa"
$d.Content.Find.Execute(
    "Questão 2 =  (primeira setença verdadeira) segunda falsa",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Questão 2 =  falsa", 2)

# Questão 3: merge two runs into "Questão 3=  falsa"
$d.Content.Find.Execute(
    "Questão 3=  ambas as setenças falsas",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Questão 3=  falsa", 2)

# Questão 6: merge many runs into "Questão 6 =  falsa"
$d.Content.Find.Execute(
    "Questão 6 =  (primeira setença falsa) (segunda setença falsa)",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Questão 6 =  falsa", 2)

# Questão 4 and Questão 5 keep the same visible text, but the diff merges
# the two runs ("Questão " + "4 = falso") into a single run. Re-run a
# Find/Replace with the same text to normalize/merge the runs.
$d.Content.Find.Execute(
    "Questão 4 = falso",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Questão 4 = falso", 2)

$d.Content.Find.Execute(
    "Questão 5 = verdadeiro",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Questão 5 = verdadeiro", 2)
